$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRows = @(2,3,5,6,7,8,9)

# ---------------------------------------------------------------------------
# 1. Capture the old values we need to relocate before we start overwriting
#    anything (old layout: C=ServiceNow ID, H=Participant Count, J=Title,
#    K=Is Active?).
# ---------------------------------------------------------------------------
$oldC = @{}
$oldH = @{}
$oldJ = @{}
$oldK = @{}
foreach ($r in $dataRows) {
    $oldC[$r] = $ws.Cells.Item($r,3).Value2
    $oldH[$r] = $ws.Cells.Item($r,8).Value2
    $oldJ[$r] = $ws.Cells.Item($r,10).Value2
    $oldK[$r] = $ws.Cells.Item($r,11).Value2
}

# ---------------------------------------------------------------------------
# 2. Re-style the header row first, while the donor cells (H1, which keeps
#    being a bold/underlined header) still carry the original style. All of
#    A1,B1,D1,E1,F1,G1,H1,I1 end up on the bold-underline header style that
#    C1/H1/J1/K1 already use.
# ---------------------------------------------------------------------------
$ws.Range("H1").Copy() | Out-Null
foreach ($c in @(1,2,4,5,6,7,9)) {
    $ws.Cells.Item(1,$c).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Propagate the "gray fill" data-style (currently on H2:H9 / K2:K9) onto
#    the new D and E data columns, cell by cell (keeps row 4 - which must
#    stay completely absent - untouched).
# ---------------------------------------------------------------------------
$ws.Range("H2").Copy() | Out-Null
foreach ($r in $dataRows) {
    $ws.Cells.Item($r,4).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

$ws.Range("K2").Copy() | Out-Null
foreach ($r in $dataRows) {
    $ws.Cells.Item($r,5).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Strip the style off the old H2:H9 cells - H becomes a plain new boolean
#    column ("Saliva Web Hooks") with no special formatting.
# ---------------------------------------------------------------------------
foreach ($r in $dataRows) {
    $ws.Cells.Item($r,8).ClearFormats() | Out-Null
}

# ---------------------------------------------------------------------------
# 5. Clear out the columns that become entirely empty in the new layout
#    (C, J, K) - removes both cell content and styling.
# ---------------------------------------------------------------------------
$ws.Range("C1:C9").Clear() | Out-Null
$ws.Range("J1:K9").Clear() | Out-Null

# ---------------------------------------------------------------------------
# 6. Header text.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Title"
$ws.Range("B1").Value = "ServiceNow ID"
$ws.Range("D1").Value = "Participant Count"
$ws.Range("E1").Value = "Is Active?"
$ws.Range("F1").Value = "Accept Saliva"
$ws.Range("G1").Value = "Accept Blood"
$ws.Range("H1").Value = "Saliva Web Hooks"
$ws.Range("I1").Value = "Blood Web Hooks"

# ---------------------------------------------------------------------------
# 7. Relocate the data values captured in step 1: A<-Title(J), B<-ServiceNow
#    ID(C), D<-Participant Count(H), E<-Is Active?(K).
# ---------------------------------------------------------------------------
foreach ($r in $dataRows) {
    $ws.Cells.Item($r,1).Value = $oldJ[$r]
    $ws.Cells.Item($r,2).Value = $oldC[$r]
    $ws.Cells.Item($r,4).Value = $oldH[$r]
    $ws.Cells.Item($r,5).Value = $oldK[$r]
}

# ---------------------------------------------------------------------------
# 8. New boolean "web hook" / "accept" columns (F,G,H,I) - plain cells, no
#    special style.
# ---------------------------------------------------------------------------
$newF = @{2=$true; 3=$true; 5=$true; 6=$true; 7=$false; 8=$false; 9=$false}
$newG = @{2=$true; 3=$true; 5=$true; 6=$false; 7=$false; 8=$false; 9=$false}
$newH = @{2=$true; 3=$true; 5=$false; 6=$false; 7=$false; 8=$false; 9=$true}
$newI = @{2=$true; 3=$false; 5=$false; 6=$false; 7=$false; 8=$true; 9=$true}
foreach ($r in $dataRows) {
    $ws.Cells.Item($r,6).Value = $newF[$r]
    $ws.Cells.Item($r,7).Value = $newG[$r]
    $ws.Cells.Item($r,8).Value = $newH[$r]
    $ws.Cells.Item($r,9).Value = $newI[$r]
}

# ---------------------------------------------------------------------------
# 9. Row 10 gains placeholder (empty) cells in F,G,H,I alongside the
#    pre-existing A10/D10. Write then clear contents while keeping the
#    "Normal" style so no explicit style id is stamped (matches A10).
# ---------------------------------------------------------------------------
foreach ($addr in @("F10","G10","H10","I10")) {
    $ws.Range($addr).Value = "x"
    $ws.Range($addr).ClearContents() | Out-Null
    $ws.Range($addr).Style = "Normal"
}

# ---------------------------------------------------------------------------
# 10. Column widths.
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 17.25      # -> stored width ~18.1666 (was 12.5 bestFit)
$ws.Columns.Item(9).ColumnWidth = 15.0       # -> stored width ~15.8333 (new column)

# ---------------------------------------------------------------------------
# 11. Selection moves from L1 to G10:I10.
# ---------------------------------------------------------------------------
$ws.Range("G10:I10").Select() | Out-Null
